# Add a new account record as row 24 of the AccountDatabase sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNum = 24

# Column order matches the sheet header: ID Number, Account Type, First Name,
# Last Name, Email, Contact Number, Nationality, Religion, Sex, Civil Status,
# Age, Disability, Permanent Address, Password.
$rowValues = @(
    "704525",
    "Facilitator",
    "dada",
    "pineda",
    "kyla@gmail.com",
    "237826",
    "filipino",
    "Catholic",
    "Male",
    "Married",
    "22",
    "no",
    "asdasd",
    "asdddd"
)

for ($i = 0; $i -lt $rowValues.Count; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item($rowNum, $col)
    # Force text storage (numeric-looking values like IDs/ages must stay
    # strings, not be coerced to numbers) without leaving a custom style
    # applied to the cell: set the number format to Text, write the value,
    # then reset the style to Normal so no style index is stamped on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $rowValues[$i]
    $cell.Style = "Normal"
}
